{"js": "const body = context.document.body;\n\n// Anchor on the text of the paragraph that precedes the block being\n// removed, rather than a hard-coded paragraph index, so the script stays\n// robust to unrelated edits earlier in the document.\nconst results = body.search(\n  \"LOT2028: Tecnologia de Processos Fermentativos (Requisito fraco)\",\n  { matchCase: true }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const anchorPara = results.items[0].paragraphs.getFirst();\n\n  // The four paragraphs immediately following the anchor are removed:\n  //   1) a blank \"Normal\" paragraph\n  //   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n  //   3) a blank \"Normal\" paragraph\n  //   4) a blank paragraph with pageBreakBefore\n  // The trailing blank paragraph and the final pageBreakBefore paragraph\n  // that follow are left untouched.\n  const p1 = anchorPara.getNext();\n  const p2 = p1.getNext();\n  const p3 = p2.getNext();\n  const p4 = p3.getNext();\n\n  p1.delete();\n  p2.delete();\n  p3.delete();\n  p4.delete();\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Requisitos\" entry paragraph that precedes the block being\n# removed (\"LOT2028: Tecnologia de Processos Fermentativos (Requisito\n# fraco)\"). Anchoring on the text (rather than a hard-coded paragraph\n# index) keeps the script robust to unrelated changes earlier in the doc.\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute(\"LOT2028: Tecnologia de Processos Fermentativos (Requisito fraco)\")\n\nif ($found) {\n    $anchorPara = $findRange.Paragraphs(1)\n\n    # The four paragraphs immediately following the anchor are removed:\n    #   1) a blank \"Normal\" paragraph\n    #   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n    #   3) a blank \"Normal\" paragraph\n    #   4) a blank paragraph with pageBreakBefore\n    # The trailing blank paragraph and the final pageBreakBefore paragraph\n    # that follow are left untouched.\n    $p1 = $anchorPara.Next()\n    $p2 = $p1.Next()\n    $p3 = $p2.Next()\n    $p4 = $p3.Next()\n\n    $deleteRange = $d.Range($p1.Range.Start, $p4.Range.End)\n    $deleteRange.Delete()\n}\n"}
